# Update countries & provincias Spain
# Applies the 13-Aug-2020 03:27 data refresh to the Pais sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp footer cell -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Agosto de 2020 a las 03:27"

# --- Rows whose country label itself changes (data re-sorted upstream) -----
$ws.Range("A170").Value = "Guadalupe"
$ws.Range("A171").Value = "Birmania"
$ws.Range("A172").Value = "Mauricio"
$ws.Range("A173").Value = "Islas Feroe"
$ws.Range("A174").Value = "Martinica"
$ws.Range("A175").Value = "Isla de Man"
$ws.Range("A176").Value = "Trinidad yTobago"

$ws.Range("A190").Value = "Polinesia Francesa"
$ws.Range("A191").Value = "Seychelles"
$ws.Range("A192").Value = "Butan"

$ws.Range("A213").Value = "Montserrat"
$ws.Range("A214").Value = "Islas Malvinas"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) -------------------

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 5360302
$ws.Range("C4").Value = 54345
$ws.Range("D4").Value = 2812576
$ws.Range("E4").Value = 2378595
$ws.Range("G4").Value = 1386
$ws.Range("H4").Value = 169131

# Row 101 - Libia
$ws.Range("B101").Value = 6611
$ws.Range("C101").Value = 309
$ws.Range("D101").Value = 778
$ws.Range("E101").Value = 5701

# Row 170 - Guadalupe
$ws.Range("B170").Value = 367
$ws.Range("C170").Value = 50
$ws.Range("D170").Value = 289
$ws.Range("E170").Value = 64
$ws.Range("H170").Value = 14

# Row 171 - Birmania
$ws.Range("B171").Value = 361
$ws.Range("C171").Value = 1
$ws.Range("D171").Value = 318
$ws.Range("E171").Value = 37
$ws.Range("H171").Value = 6

# Row 172 - Mauricio
$ws.Range("B172").Value = 344
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 334
$ws.Range("E172").Value = 0
$ws.Range("H172").Value = 10

# Row 173 - Islas Feroe
$ws.Range("B173").Value = 339
$ws.Range("C173").Value = 21
$ws.Range("D173").Value = 225
$ws.Range("E173").Value = 114
$ws.Range("H173").Value = 0

# Row 174 - Martinica
$ws.Range("D174").Value = 98
$ws.Range("E174").Value = 222
$ws.Range("H174").Value = 16

# Row 175 - Isla de Man
$ws.Range("B175").Value = 336
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 312
$ws.Range("E175").Value = 0
$ws.Range("H175").Value = 24

# Row 176 - Trinidad yTobago
$ws.Range("B176").Value = 326
$ws.Range("C176").Value = 26
$ws.Range("D176").Value = 139
$ws.Range("E176").Value = 179
$ws.Range("H176").Value = 8

# Row 187 - Barbados
$ws.Range("B187").Value = 144
$ws.Range("C187").Value = 1
$ws.Range("D187").Value = 115
$ws.Range("E187").Value = 22

# Row 190 - Polinesia Francesa
$ws.Range("B190").Value = 139
$ws.Range("C190").Value = 27
$ws.Range("D190").Value = 64
$ws.Range("E190").Value = 75

# Row 191 - Seychelles
$ws.Range("B191").Value = 127
$ws.Range("D191").Value = 126
$ws.Range("E191").Value = 1

# Row 192 - Butan
$ws.Range("B192").Value = 113
$ws.Range("D192").Value = 97
$ws.Range("E192").Value = 16

# Row 213 - Montserrat
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

# Row 214 - Islas Malvinas
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
